# Generate Report for Handback
#
# The localization run finished handing back `d74854b4-...md` (it is now
# "in sync with en-US" for both target locales), so the status report is
# regenerated: the "Ready for handoff" rows flip to
# "Handed back: in sync with en-US", the handback timestamps advance, and
# the stale "version mismatch" error detail is cleared now that the
# handback is current.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"

# ---- zh-cn sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K3").Value = "2016-08-31 18:56:09"
$ws.Range("P3").Value = ""
$ws.Columns.Item(16).ColumnWidth = 13.7470528738839

# ---- de-de sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("K3").Value = "2016-08-31 18:56:18"
$ws.Range("P3").Value = ""
$ws.Columns.Item(16).ColumnWidth = 13.7470528738839
